$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# H2 held " A1" (leading space) as a stray/duplicate value; normalize it to "A1".
$ws.Range("H2").Value = "A1"

# D2's dvno changes from 7080 to 8010 (testing bill csv data swapped in).
$ws.Range("D2").Value = 8010

# Rows 3-6 were extra duplicate bill rows (dvno 7021-7024, bundleno "A4");
# they are removed entirely, shrinking the sheet down to just the header + one data row.
$ws.Range("A3:I6").EntireRow.Delete()
